$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = 2197704.26
$ws.Range("C7").Value = -51.35059287224324
$ws.Range("D7").Value = 2084
$ws.Range("E7").Value = 2084
$ws.Range("F7").Value = 1054.560585412668
$ws.Range("G7").Value = 8.830871415356079
